$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "ema" header in G1, copying the bold/centered/bordered
# header style from the existing F1 ("Volume") header cell.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "ema"

# EMA(7) of the Close column (E), seeded with the first close value.
$emaValues = @(
3769.697021484375,
3784.6640014648438,
3778.8430938720703,
3782.6464805603027,
3724.5816011428833,
3648.0382516384125,
3534.331301033497,
3423.7416398376226,
3357.2440960891545,
3288.707547164522,
3276.058553439798,
3300.108490275161,
3287.153511261058,
3292.865499656731,
3302.281815172236,
3314.441830129177,
3288.907605507039,
3257.686990751373,
3217.2217249971236,
3163.1963230447177,
3011.880152439788,
2860.2054024157787,
2778.913817436834,
2694.273436808094,
2634.688842254508,
2593.024200050256,
2550.5184552134733,
2549.661851664011,
2561.517567947227,
2572.004814144014,
2601.073312756448,
2648.834281442336,
2657.3392120583144,
2662.7950584577984,
2742.9930223589736,
2810.906815109074,
2872.549129886493,
2940.02952319612,
2985.674303041621,
3049.1199850937155,
3056.2105332538804,
3024.003786170879,
2997.3435378703466,
2968.8734981293223,
2960.0248745735544,
3014.987955246572,
3043.198485966179,
3002.7693478730716,
2948.5088956704285,
2902.3069647215716,
2833.892332916179,
2768.873290214478,
2736.479796762421,
2699.9497523569717,
2674.4790989356975,
2696.9932353345857,
2718.0228806025016,
2693.9675999050014,
2750.275992897501,
2805.8282715285945,
2841.900805697227,
2840.0428528080765,
2784.321141070901,
2754.4486194750507,
2704.5958029851945,
2652.88966229749,
2633.854136371555,
2657.836464095072,
2645.389494067398,
2623.932867620861,
2611.638188313302,
2588.464786254508,
2589.0226009213497,
2596.8043652417937,
2640.61718994697,
2684.1765155071025,
2749.468141024858,
2798.665375788175,
2814.1138389700372,
2835.079519852528,
2869.592415768302,
2909.9610964941953,
2959.486325300334,
2996.282590654938,
3033.006681760735,
3097.6494205002386,
3157.395695746273,
3218.5437200518923,
3260.197279785013,
3265.558665405166,
3311.557060577312,
3344.9326269759526,
3389.407844255402,
3422.3661859259264,
3419.7227595616323,
3357.715043304037,
3326.604947028809,
3292.9722039239505,
3285.2082178843693,
3266.8728992531205,
3195.4177359632777,
3154.1574181833957,
3145.204123207859,
3113.8804483629256,
3095.639462248757,
3087.3071723701614,
3063.831209355746,
3062.2750793800906,
3072.7329208631927,
3073.9861530497383,
3052.35979300996,
3030.4787680973136,
3007.387579490954,
2986.2238511221217,
2992.016277013466,
2946.086792721037,
2931.797516415778,
2933.083366804021,
2903.712976763172,
2860.33142446691,
2852.187594229089,
2853.4932957694728,
2835.9890514169483,
2862.1529824103673,
2833.9180204991817,
2799.1834494564173,
2758.410841486844,
2698.173121349508,
2584.9874459926,
2524.6183310764814,
2411.4909089518924,
2299.043566235404,
2227.8872278992094,
2184.9839023697195,
2175.164623554633,
2137.0549557030845,
2125.3935116991884,
2073.209165756813,
2059.490919727766,
2034.9471082528557,
2019.8399088263604,
2025.6724731236766,
2012.2998260341637,
2003.9705665471074,
1989.1848951251743,
1942.867003863412,
1888.380965788184,
1855.7711674563725,
1844.8361270571231,
1882.7374163377644,
1897.6350629857452,
1879.118631223684,
1867.8766015915908,
1844.6771045140056,
1833.9102026530823,
1826.7338910034837,
1834.8728364655035,
1829.6667123100651,
1820.6431006387988,
1812.938837930271,
1775.9646875297344,
1714.38937868441,
1647.0961721969013,
1536.4678200656447,
1455.2665754984523,
1399.7515368484487,
1316.7463308589927,
1259.1895699215884,
1192.8013723874803,
1176.5151528257666,
1164.2969786329968,
1154.4288802149822,
1128.6771289112367,
1132.3545263709275,
1155.9770764188206,
1177.8445321676309,
1183.3413153854885,
1185.9261525547413,
1175.589420324259,
1156.4280271572568,
1134.1457273991925,
1115.551129045488,
1103.2915511298193,
1095.910405290724,
1109.6975744758554,
1115.9084347631415,
1133.6748258282155,
1159.6544653184271,
1175.3674053853047,
1185.7701219100722,
1181.4279942645853,
1160.3801387648452,
1129.8330166712901,
1125.7715520542488,
1142.2102253199837,
1164.939956587644,
1211.8615836516706,
1243.5551232856278,
1327.3458163411738,
1381.2531671875208,
1415.9900462890782,
1456.1799126464962,
1476.4862223266691,
1494.6890380828922,
1520.8859497047474,
1502.0103179914513,
1486.9594291674166,
1524.2777383794687,
1574.5753386967108,
1612.7832496280018,
1633.5798078264702,
1645.5641893659463,
1642.972092219772,
1640.4654278074072,
1635.067698785243,
1628.352226725651,
1654.327829712207,
1663.660392547827,
1672.5830019304017,
1698.3162797681139,
1699.493465929601,
1737.555768392513,
1773.4728565678224,
1819.4162574161014,
1859.8963300249666,
1879.122748007006,
1885.39909835877,
1883.5841748432963,
1870.938039579738,
1864.9554828098035,
1801.9634382792276,
1745.7235247543426,
1714.1223859973975,
1691.218254341798,
1684.1061651215828,
1677.3444248665778,
1682.1225764624332,
1638.5376415753403,
1601.7519860643179,
1558.950830368551,
1557.4724611553195,
1549.0640626633647,
1550.219281128383,
1559.2086502994123,
1563.711602470653,
1562.001872019005,
1565.9118044048787,
1578.7296601884245,
1574.484378930381,
1588.3398772153637,
1600.0918219740229,
1629.8402287852045,
1666.4311054267941,
1690.2733412771267,
1696.1463206551107,
1667.30673145813,
1659.1687998143007,
1612.2999702220536,
1567.3369161431028,
1542.9381119120146,
1491.035859324636,
1462.6622399524613,
1428.0937258627835,
1384.2222414185721,
1370.0867250092415,
1359.6299241280249,
1349.2207646292218,
1335.4697726906663,
1335.4323771254217,
1334.106204230785,
1334.9323753410574,
1335.1123674432931,
1333.3289350062978,
1327.9078035496455,
1314.9542291270777,
1317.0754740914022,
1328.3382522970671,
1334.4629787735817,
1338.7746022442489,
1337.2101936265462,
1331.7827672902222,
1329.4881374793854,
1319.9505147306327,
1309.856806946412,
1306.1191982273872,
1301.6203740123374,
1300.5708102944093,
1294.1460337461976,
1297.183683024492,
1305.8161680300877,
1306.9738813936597,
1301.6664767190728,
1297.0500894728984,
1297.7741698878767,
1301.905425999892,
1317.2908248710128,
1324.217752442322,
1358.579665650101,
1410.5764118352322,
1436.526028358846,
1466.263997587494,
1504.622619772652,
1526.162796128317,
1537.8007164810033,
1548.27668482169,
1541.1354616143144,
1538.7370332224546,
1565.3261208641065,
1580.986595042611,
1578.7986315846929,
1576.2468008369572,
1515.3939934499836,
1411.5879450386597,
1383.5571086813384,
1359.4730988449883,
1333.42190177046,
1305.5212297946418,
1289.5419843577001,
1280.090539781947,
1263.9685396020852,
1248.178553139064,
1239.2089880964854,
1234.0134305254892,
1211.1267415952889,
1185.4333130324042,
1172.8683502528188,
1175.4511589298486,
1182.584157771605,
1186.6695941587818,
1191.4766829237738,
1192.3892504740802,
1186.8134832657165,
1194.3354237285844,
1219.6737174546413,
1233.8237695362934,
1248.9436633338605,
1247.5414572171924,
1255.720240862113,
1256.7093700997098,
1260.4454797232197,
1253.443484792415,
1260.3616968462643,
1261.3424701444637,
1262.6028914267072,
1262.91930113839,
1265.8442304924645,
1279.520471453333,
1286.97253742789,
1281.817873529902,
1253.428255000942,
1237.1085752839094,
1224.0102339531666,
1209.9101412851874,
1211.8585092842031,
1212.293875859637,
1213.76593912129,
1215.3643102979988,
1216.810372616077,
1217.348288495261,
1219.7548076800394,
1218.0140122541702,
1211.0070301867213,
1208.6541068685565,
1206.2987771729017
)

for ($i = 0; $i -lt $emaValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $emaValues[$i]
}
